$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "Enquiry" columns for the CurrencyRates param (Z & AA), rows 2-3.
$ws.Range("Z2").Value = "For Multi Currency Transactions, the Exchange Rates (Buy and Sell)"
$ws.Range("AA2").Value = "Choice of various modes of payments "
$ws.Range("Z3").Value = "Rate of Buy / Sell "

# Apply the formatting (7.5pt font, thin box border, wrap + vertically centered)
# on Z2 first, then replicate it onto AA2 / Z3 via a format-only paste so the
# workbook ends up with a single shared cell style instead of one per cell.
$src = $ws.Range("Z2")
$src.Borders.LineStyle = 1
$src.Font.Size = 7.5
$src.VerticalAlignment = -4108
$src.WrapText = $true

$src.Copy()
$ws.Range("AA2").PasteSpecial(-4122)
$ws.Range("Z3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Widen the two new columns to fit the longer descriptions.
$ws.Columns.Item(26).ColumnWidth = 53
$ws.Columns.Item(27).ColumnWidth = 52

# Match the saved selection/viewport from the authored workbook.
$ws.Range("Z2:AA3").Select()
